$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$grp = $s.Shapes.Item(2)

# --- Shape "TextBox 25" (id 26): reword "down key" -> "right key" and
# nudge its position slightly, splitting the sentence across 3 runs
# (same formatting on every run) exactly like the authored diff.
$sh25 = $grp.GroupItems.Item(24)

$tr25 = $sh25.TextFrame.TextRange
$sub25 = $tr25.Characters(1, 40)
$sub25.Text = "1. Press the arrow right key on a focused"
# Touching the colour of the middle word (even though it is re-set to the
# very same colour) is what forces the host to persist the text as three
# separate <a:r> runs instead of collapsing them back into one.
$mid25 = $tr25.Characters(20, 10)
$mid25.Font.Color.RGB = 0x0000C0

$sh25.Left = 116.88220977783203
$sh25.Top = 66.24024200439453

# --- Shape "TextBox 34" (id 35): position nudge only.
$sh34 = $grp.GroupItems.Item(26)
$sh34.Left = 496.512939453125
$sh34.Top = 66.6912612915039

# --- Shape "TextBox 61" (id 62): position nudge only.
$sh61 = $grp.GroupItems.Item(28)
$sh61.Left = 116.88220977783203
$sh61.Top = 313.5706481933594

# --- Shape "TextBox 63" (id 64): position nudge only.
$sh63 = $grp.GroupItems.Item(30)
$sh63.Left = 496.85198974609375
$sh63.Top = 307.70538330078125
